$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer IDs for the Deleted Customer test case (rows 13-18, column A)
$ids = @(
    "cus_IBc0ERhRyxXWsL",
    "cus_IBc00mEJZdW8Kg",
    "cus_IBc09g64O3FaQE",
    "cus_IBc0ej42CA1Txb",
    "cus_IBc0J46XloVal7",
    "cus_IBc0swmV0KXVgB"
)

$row = 13
foreach ($id in $ids) {
    $ws.Cells.Item($row, 1).Value = $id
    $row = $row + 1
}

# Update the selection to match the authored state
$ws.Range("H15").Select()
